# Update the 진단명 (diagnosis) column values to the new "질병명_N" naming
# scheme (previously "질병_N").
#
# NOTE on write order: the workbook's shared-string table appends newly
# introduced strings in the order cells are actually written (and drops
# strings that become unreferenced). To reproduce the exact shared-string
# layout of the target file, row 6 is written before row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "질병명_1"
$ws.Range("G3").Value = "질병명_1"
$ws.Range("G4").Value = "질병명_2,질병명_4"
$ws.Range("G6").Value = "질병명_7"
$ws.Range("G5").Value = "질병명_1,질병명_3"
$ws.Range("G7").Value = "질병명_3,질병명_8"
$ws.Range("G8").Value = "질병명_2,질병명_3"

# Move the saved selection to match the author's last on-screen state.
$ws.Range("G9").Select()
